$wb = $excel.ActiveWorkbook

# Sheet ALC, row 3
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 42925.6
$ws.Range("J3").Value = 42925.6
$ws.Range("L3").Value = 42925.6
$ws.Range("N3").Value = -43153.6

# Sheet ALC, row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3075.25
$ws.Range("I18").Value = 1912.375
$ws.Range("J18").Value = 5401
$ws.Range("K18").Value = 1912.375
$ws.Range("L18").Value = 5401
$ws.Range("M18").Value = -1628.375
$ws.Range("N18").Value = -5969

# Sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 347.04166
$ws.Range("I33").Value = 376.95
$ws.Range("J33").Value = 197.5
$ws.Range("K33").Value = 376.95
$ws.Range("L33").Value = 197.5
$ws.Range("M33").Value = -147.95
$ws.Range("N33").Value = -655.5

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1686.1111
$ws.Range("I40").Value = 1827.7778
$ws.Range("J40").Value = 1544.4445
$ws.Range("K40").Value = 1827.7778
$ws.Range("L40").Value = 1544.4445
$ws.Range("M40").Value = -1652.7778
$ws.Range("N40").Value = -1894.4445

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4185.6562
$ws.Range("I43").Value = 2311.4285
$ws.Range("J43").Value = 4710.44
$ws.Range("K43").Value = 2311.4285
$ws.Range("L43").Value = 4710.44
$ws.Range("M43").Value = -2242.4285
$ws.Range("N43").Value = -4848.44

# Sheet ALC, row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5131708
$ws.Range("I64").Value = 7695652
$ws.Range("J64").Value = 3820
$ws.Range("K64").Value = 7695652
$ws.Range("L64").Value = 3820
$ws.Range("M64").Value = -7695404
$ws.Range("N64").Value = -4316

# Sheet ALC, row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5131708
$ws.Range("I67").Value = 7695652
$ws.Range("J67").Value = 3820
$ws.Range("K67").Value = 7695652
$ws.Range("L67").Value = 3820
$ws.Range("M67").Value = -7694794
$ws.Range("N67").Value = -5536

# Sheet ALC, row 102
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 42925.6
$ws.Range("J102").Value = 42925.6
$ws.Range("L102").Value = 42925.6
$ws.Range("N102").Value = -49415.6

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8334373.5
$ws.Range("I137").Value = 806.6667
$ws.Range("J137").Value = 11112229
$ws.Range("K137").Value = 2420.0001
$ws.Range("L137").Value = 33336687
$ws.Range("M137").Value = 129.9998999999998
$ws.Range("N137").Value = -33341787

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5651316
$ws.Range("I138").Value = 7577070
$ws.Range("J138").Value = 2436.6667
$ws.Range("K138").Value = 22731210
$ws.Range("L138").Value = 7310.000100000001
$ws.Range("M138").Value = -22726070
$ws.Range("N138").Value = -17590.0001

# Sheet ALC, row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 49720
$ws.Range("J140").Value = 49720
$ws.Range("L140").Value = 49720
$ws.Range("N140").Value = -60080

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10599.288
$ws.Range("I32").Value = 11570.945
$ws.Range("J32").Value = 7630.3335
$ws.Range("K32").Value = 11570.945
$ws.Range("L32").Value = 7630.3335
$ws.Range("M32").Value = -11283.945
$ws.Range("N32").Value = -8204.333500000001

# Sheet ARM, row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2462.5
$ws.Range("I63").Value = 2426.875
$ws.Range("J63").Value = 2510
$ws.Range("K63").Value = 2426.875
$ws.Range("L63").Value = 2510
$ws.Range("M63").Value = -1740.875
$ws.Range("N63").Value = -3882

# Sheet ARM, row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2462.5
$ws.Range("I66").Value = 2426.875
$ws.Range("J66").Value = 2510
$ws.Range("K66").Value = 12134.375
$ws.Range("L66").Value = 12550
$ws.Range("M66").Value = -8702.375
$ws.Range("N66").Value = -19414

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9435555
$ws.Range("I74").Value = 12501557
$ws.Range("J74").Value = 1701.6923
$ws.Range("K74").Value = 12501557
$ws.Range("L74").Value = 1701.6923
$ws.Range("M74").Value = -12500683
$ws.Range("N74").Value = -3449.6923

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9435555
$ws.Range("I77").Value = 12501557
$ws.Range("J77").Value = 1701.6923
$ws.Range("K77").Value = 62507785
$ws.Range("L77").Value = 8508.461499999999
$ws.Range("M77").Value = -62503417
$ws.Range("N77").Value = -17244.4615

# Sheet BSM, row 5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 523.6
$ws.Range("I5").Value = 404.5
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 404.5
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -291.5
$ws.Range("N5").Value = -1226

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 27780024
$ws.Range("I86").Value = 1900
$ws.Range("K86").Value = 1900
$ws.Range("M86").Value = -777

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 27780024
$ws.Range("I89").Value = 1900
$ws.Range("K89").Value = 9500
$ws.Range("M89").Value = -3884

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2174.8833
$ws.Range("I134").Value = 1327.7906
$ws.Range("K134").Value = 3983.3718
$ws.Range("M134").Value = -1448.3718

# Sheet CRP, row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 14334.875
$ws.Range("I59").Value = 5104
$ws.Range("J59").Value = 15653.571
$ws.Range("K59").Value = 5104
$ws.Range("L59").Value = 15653.571
$ws.Range("M59").Value = -3959
$ws.Range("N59").Value = -17943.571

# Sheet CRP, row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 8030.6
$ws.Range("J60").Value = 8551
$ws.Range("L60").Value = 8551
$ws.Range("N60").Value = -9573

# Sheet CRP, row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2247.4119
$ws.Range("I62").Value = 2203.5715
$ws.Range("K62").Value = 2203.5715
$ws.Range("M62").Value = -1579.5715

# Sheet CRP, row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2247.4119
$ws.Range("I65").Value = 2203.5715
$ws.Range("K65").Value = 11017.8575
$ws.Range("M65").Value = -7897.8575

# Sheet CRP, row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17799
$ws.Range("J68").Value = 17799
$ws.Range("L68").Value = 17799
$ws.Range("N68").Value = -19297

# Sheet CRP, row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17799
$ws.Range("J71").Value = 17799
$ws.Range("L71").Value = 53397
$ws.Range("N71").Value = -60885

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 760
$ws.Range("I107").Value = 751.25
$ws.Range("J107").Value = 795
$ws.Range("K107").Value = 751.25
$ws.Range("L107").Value = 795
$ws.Range("M107").Value = 1168.75
$ws.Range("N107").Value = -4635

# Sheet CRP, row 120
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 39071.855
$ws.Range("J120").Value = 39071.855
$ws.Range("L120").Value = 39071.855
$ws.Range("N120").Value = -46329.855

# Sheet CUL, row 42
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 1899.75
$ws.Range("J42").Value = 2199.6667
$ws.Range("L42").Value = 6599.000100000001
$ws.Range("N42").Value = -7667.000100000001

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4991.1914
$ws.Range("I7").Value = 4990.231
$ws.Range("J7").Value = 4992.381
$ws.Range("K7").Value = 4990.231
$ws.Range("L7").Value = 4992.381
$ws.Range("M7").Value = -4878.231
$ws.Range("N7").Value = -5216.381

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 878.6
$ws.Range("I46").Value = 574.2857
$ws.Range("J46").Value = 1144.875
$ws.Range("K46").Value = 574.2857
$ws.Range("L46").Value = 1144.875
$ws.Range("M46").Value = -386.2857
$ws.Range("N46").Value = -1520.875

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1529.375
$ws.Range("I68").Value = 1787
$ws.Range("J68").Value = 1100
$ws.Range("K68").Value = 1787
$ws.Range("L68").Value = 1100
$ws.Range("M68").Value = -1038
$ws.Range("N68").Value = -2598

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1529.375
$ws.Range("I71").Value = 1787
$ws.Range("J71").Value = 1100
$ws.Range("K71").Value = 8935
$ws.Range("L71").Value = 5500
$ws.Range("M71").Value = -5191
$ws.Range("N71").Value = -12988

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4991.1914
$ws.Range("I126").Value = 4990.231
$ws.Range("J126").Value = 4992.381
$ws.Range("K126").Value = 14970.693
$ws.Range("L126").Value = 14977.143
$ws.Range("M126").Value = -12500.693
$ws.Range("N126").Value = -19917.143

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6950372.5
$ws.Range("I132").Value = 3585.8823
$ws.Range("J132").Value = 23821140
$ws.Range("K132").Value = 10757.6469
$ws.Range("L132").Value = 71463420
$ws.Range("M132").Value = -8227.6469
$ws.Range("N132").Value = -71468480

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 629.7241
$ws.Range("I107").Value = 709.5217
$ws.Range("J107").Value = 323.83334
$ws.Range("K107").Value = 2128.5651
$ws.Range("L107").Value = 971.5000200000001
$ws.Range("M107").Value = -208.5650999999998
$ws.Range("N107").Value = -4811.50002
